# Pull request to merge ui and api
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# EXECUTE flag for TC_002_AMAZON_LOGIN (row 3) switches from YES to NO
$ws.Range("B3").Value = "NO"

# PARALLEL flag for the last existing UI row (row 15) now set explicitly to NO
$ws.Range("H15").Value = "NO"

# Row 16 previously held only a tall, placeholder row format with no data.
# Reset it back to a normal row/style before filling it with the new API
# test-case data so it looks like a freshly authored row.
$ws.Range("C16").Style = "Normal"
$ws.Rows("16").AutoFit()

# New API test rows appended below the UI rows
$apiRows = @(
    @("TC_014_GET_CALL", "YES", "API", "ApiCalls", "getSingleUserCall",   "data.xlsx", "API_DATA_FILE", "NO"),
    @("TC_014_GET_CALL", "YES", "API", "ApiCalls", "postCreateUserCall",  "data.xlsx", "API_DATA_FILE", "NO"),
    @("TC_014_GET_CALL", "YES", "API", "ApiCalls", "putUpdateUserCall",   "data.xlsx", "API_DATA_FILE", "NO"),
    @("TC_014_GET_CALL", "YES", "API", "ApiCalls", "deleteUserCall",      "data.xlsx", "API_DATA_FILE", "NO")
)

$startRow = 16
for ($i = 0; $i -lt $apiRows.Count; $i++) {
    $r = $startRow + $i
    $row = $apiRows[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

# Column H (PARALLEL) on the brand-new rows is filled in by carrying the
# formatting down from the first new row so it keeps the same look as the
# rest of the sheet.
$ws.Range("H16").Copy()
$ws.Range("H17:H19").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Scroll the window so row 7 is at the top and select C16 (where the author
# was working when the new API rows were added).
$win = $wb.Windows.Item(1)
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("C16").Select()
